$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('account', '口座|こうざ')
    ,@('commission', '手数料|てすうりょう')
    ,@('bank card', 'キャッシュカード')
    ,@('passbook', '通帳|つうちょう')
    ,@('savings', '預金|よきん')
    ,@('amount', '金額|きんがく')
    ,@('personal identification number', '暗証番号|あんしょうばんごう')
    ,@('cash', '現金|げんきん')
    ,@('money exchange', '両替|りょうがえ')
    ,@('deposit', 'お預け入れ|おあずけいれ')
    ,@('withdrawal', 'お引き出し|おひきだし')
    ,@('passbook update', '通帳記入|つうちょうきにゅう')
    ,@('bank transfer', 'お振込|おふりこみ')
    ,@('credit card', 'クレジットカード')
    ,@('balance inquiry', '残高照会|ざんだかしょうかい')
    ,@('I would like to open an account.', '口座を開きたいんですが。|こうざをひらきたいんですが。')
    ,@('I would like to close an account.', '口座を閉じたいんですが。|こうざをとじたいんですが。')
    ,@('Please change dollars into yen.', 'ドルを円にかえてください。|ドルを円にかえてください。')
    ,@('I would like to deposit money into the account.', '口座にお金を振り込みたいんですが。|こうざにおかねをふりこみたいんですが。')
    ,@('Can you change a 10,000-yen bill into ten 1,000-yen bills?', '一万円札を千円札十枚に両替できますか。|いちまんえんさつをせんえんさつじゅうまいにりょうがえできますか。')
    ,@('I will withdraw money.', 'お金をおろします。|おかねをおろします。')
)

$startRow = 58
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
